$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 (I0) and J1 (IF), matching style of existing header (H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill in I and J values for rows 2-84 (parallel arrays, this runtime flattens nested arrays)
$rowNums = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50, 51, 52, 53, 54, 55, 56, 57, 58, 59, 60, 61, 62, 63, 64, 65, 66, 67, 68, 69, 70, 71, 72, 73, 74, 75, 76, 77, 78, 79, 80, 81, 82, 83, 84)
$iVals = @(7, 6, 7, 6, 6, 7, 7, 6, 8, 9, 10, 7, 9, 8, 8, 7, 9, 10, 8, 10, 7, 9, 8, 8, 8, 10, 8, 8, 7, 9, 8, 8, 8, 10, 8, 9, 9, 8, 7, 8, 9, 7, 7, 7, 8, 7, 7, 8, 9, 8, 8, 9, 9, 8, 8, 8, 7, 8, 8, 8, 8, 8, 9, 9, 7, 8, 8, 9, 8, 9, 9, 9, 9, 9, 9, 8, 8, 8, 8, 4, 7, 4, 4)
$jVals = @(7, 6, 7, 6, 6, 7, 7, 6, 8, 9, 10, 7, 9, 8, 8, 7, 9, 10, 8, 10, 7, 9, 8, 8, 9, 10, 8, 8, 7, 9, 8, 8, 8, 10, 8, 9, 9, 8, 7, 8, 9, 7, 7, 7, 8, 7, 7, 8, 9, 8, 8, 9, 9, 8, 8, 8, 7, 8, 8, 8, 8, 8, 9, 9, 7, 8, 8, 9, 8, 9, 9, 9, 9, 9, 9, 8, 8, 8, 8, 4, 7, 4, 4)

for ($k = 0; $k -lt $rowNums.Count; $k++) {
    $ws.Cells.Item($rowNums[$k], 9).Value = $iVals[$k]
    $ws.Cells.Item($rowNums[$k], 10).Value = $jVals[$k]
}
